# "Generate Report for handback" - refresh the handback status report.
#
# Both localized files (26508feb... and 78d9e4aa...) have now been handed
# back, so every sheet is regenerated with:
#   - rows re-sorted alphabetically by source file name
#     (26508feb-... now sorts before 78d9e4aa-...)
#   - status flipped from "Ready for handoff" to "Handed back: in sync with en-US"
#   - a fresh "Latest Handback DateTime" stamp for both locales

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws.Range("B2").Value = $statusHandedBack
$ws.Range("C2").Value = $statusHandedBack

$ws.Range("A3").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws.Range("B3").Value = $statusHandedBack
$ws.Range("C3").Value = $statusHandedBack

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/37aa0ba8dc7d88581e611cb910bbbc13553b7e43/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/37aa0ba8dc7d88581e611cb910bbbc13553b7e43/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/37aa0ba8dc7d88581e611cb910bbbc13553b7e43/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws.Range("B2").Value = $statusHandedBack
$ws.Range("C2").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf"
$ws.Range("D2").Value = "2016-01-26 12:30:31"
$ws.Range("E2").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws.Range("F2").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf"
$ws.Range("G2").Value = "2016-01-26 12:31:21"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws.Range("B3").Value = $statusHandedBack
$ws.Range("C3").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf"
$ws.Range("D3").Value = "2016-01-26 12:30:31"
$ws.Range("E3").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws.Range("F3").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf"
$ws.Range("G3").Value = "2016-01-26 12:31:21"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/37aa0ba8dc7d88581e611cb910bbbc13553b7e43/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d56c184cf810521cfee8d602c5bdb8b94110d7e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/05fef02f6be56dc738d88cad2136542ecec9edcc/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/601b964838776b0b272433c2c758816aa5a94328/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/37aa0ba8dc7d88581e611cb910bbbc13553b7e43/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d56c184cf810521cfee8d602c5bdb8b94110d7e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/05fef02f6be56dc738d88cad2136542ecec9edcc/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/601b964838776b0b272433c2c758816aa5a94328/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/37aa0ba8dc7d88581e611cb910bbbc13553b7e43/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws.Range("B2").Value = $statusHandedBack
$ws.Range("C2").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf"
$ws.Range("D2").Value = "2016-01-26 12:30:43"
$ws.Range("E2").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md"
$ws.Range("F2").Value = "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf"
$ws.Range("G2").Value = "2016-01-26 12:31:41"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws.Range("B3").Value = $statusHandedBack
$ws.Range("C3").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf"
$ws.Range("D3").Value = "2016-01-26 12:30:43"
$ws.Range("E3").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md"
$ws.Range("F3").Value = "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf"
$ws.Range("G3").Value = "2016-01-26 12:31:41"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/37aa0ba8dc7d88581e611cb910bbbc13553b7e43/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f794acf7fade616765620eb6848af0c337354e6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e67c4144042519038dd518a6b18fdab830a05744/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b29d735140eb965d8d21e791dcf322346526a9ac/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/37aa0ba8dc7d88581e611cb910bbbc13553b7e43/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f794acf7fade616765620eb6848af0c337354e6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e67c4144042519038dd518a6b18fdab830a05744/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b29d735140eb965d8d21e791dcf322346526a9ac/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/37aa0ba8dc7d88581e611cb910bbbc13553b7e43/.localization-config", "", "", ".localization-config")

Write-Output "done"
